$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "66.982.79"
$ws.Cells.Item(2, 5).Value = "  -3.71%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.380.03"
$ws.Cells.Item(3, 5).Value = "  -4.62%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "562.77"
$ws.Cells.Item(5, 5).Value = "  -4.04%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "185.41"
$ws.Cells.Item(6, 5).Value = "  -6.22%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.601"
$ws.Cells.Item(7, 5).Value = "  -2.09%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.06%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "3.370.97"
$ws.Cells.Item(9, 5).Value = "  -4.42%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.190"
$ws.Cells.Item(10, 5).Value = "  -8.80%  "
$ws.Cells.Item(11, 5).Value = "  -4.84%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "48.24"
$ws.Cells.Item(12, 5).Value = "  -7.45%  "
$ws.Cells.Item(13, 5).Value = "  -5.96%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "8.80"
$ws.Cells.Item(14, 5).Value = "  -6.02%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.912.54"
$ws.Cells.Item(15, 5).Value = "  -4.47%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "607.62"
$ws.Cells.Item(16, 5).Value = "  -10.67%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "66.797.88"
$ws.Cells.Item(17, 5).Value = "  -3.94%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.375.95"
$ws.Cells.Item(18, 5).Value = "  -4.03%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "18.06"
$ws.Cells.Item(19, 5).Value = "  -3.08%  "
$ws.Cells.Item(20, 5).Value = "  -2.83%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.73"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.920"
$ws.Cells.Item(22, 5).Value = "  -5.32%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "17.04"
$ws.Cells.Item(23, 5).Value = "  -5.07%  "
$ws.Cells.Item(24, 5).Value = "  -2.28%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "98.93"
$ws.Cells.Item(25, 5).Value = "  -8.24%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "4.12"
$ws.Cells.Item(26, 5).Value = "  -6.90%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.02"
$ws.Cells.Item(27, 5).Value = "  +0.40%  "
$ws.Cells.Item(28, 5).Value = "  -6.39%  "
$ws.Cells.Item(29, 5).Value = "  -7.80%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.88"
$ws.Cells.Item(30, 5).Value = "  -8.83%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "31.00"
$ws.Cells.Item(31, 5).Value = "  -7.47%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.98"
$ws.Cells.Item(32, 5).Value = "  -9.75%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.38"
$ws.Cells.Item(33, 5).Value = "  -8.13%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "11.25"
$ws.Cells.Item(34, 5).Value = "  -5.94%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "554.18"
$ws.Cells.Item(35, 5).Value = "  +9.61%  "
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.107"
$ws.Cells.Item(36, 5).Value = "  -5.07%  "
$ws.Cells.Item(37, 2).Value = "Maker"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.867.64"
$ws.Cells.Item(37, 5).Value = "  +1.66%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "58.78"
$ws.Cells.Item(38, 5).Value = "  -5.69%  "
$ws.Cells.Item(39, 5).Value = "  -0.03%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.45"
$ws.Cells.Item(40, 5).Value = "  -4.74%  "
$ws.Cells.Item(41, 2).Value = "CoreDAO"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.52"
$ws.Cells.Item(41, 5).Value = "  +29.95%  "
$ws.Cells.Item(42, 2).Value = "PEPE"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0₃0729"
$ws.Cells.Item(42, 5).Value = "  -11.24%  "
$ws.Cells.Item(43, 2).Value = "Fetch.AI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.74"
$ws.Cells.Item(43, 5).Value = "  -8.08%  "
$ws.Cells.Item(44, 2).Value = "Kaspa"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.129"
$ws.Cells.Item(44, 5).Value = "  -5.36%  "
$ws.Cells.Item(45, 5).Value = "  -5.73%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "32.59"
$ws.Cells.Item(46, 5).Value = "  -6.86%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0423"
$ws.Cells.Item(47, 5).Value = "  -8.26%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.26"
$ws.Cells.Item(48, 5).Value = "  -3.47%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.70"
$ws.Cells.Item(49, 5).Value = "  -9.29%  "
$ws.Cells.Item(50, 5).Value = "  -4.98%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.998"
$ws.Cells.Item(51, 5).Value = "  -0.09%  "
